# New crime data collected - update the 111th Precinct weekly CompStat report
# for the week of 1/15/2024 - 1/21/2024 (Volume 31, Number 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (edit in place so unrelated rich-text runs are left
# alone as much as possible).
# ---------------------------------------------------------------------------

# A8: "Volume 31   Number  2" -> "Volume 31   Number  3"
$a8 = $ws.Range("A8")
$a8Text = $a8.Value2
$a8Idx = $a8Text.LastIndexOf("2") + 1
$a8.Characters($a8Idx, 1).Text = "3"

# C9: "Report Covering the Week  1/8/2024  Through  1/14/2024"
#  -> "Report Covering the Week  1/15/2024  Through  1/21/2024"
$c9 = $ws.Range("C9")
$c9Text = $c9.Value2
$c9Idx1 = $c9Text.IndexOf("1/8/2024") + 1
$c9.Characters($c9Idx1, 8).Text = "1/15/2024"
$c9Text2 = $c9.Value2
$c9Idx2 = $c9Text2.IndexOf("1/14/2024") + 1
$c9.Characters($c9Idx2, 9).Text = "1/21/2024"

# ---------------------------------------------------------------------------
# Helper functions for writing the crime-stat table while keeping the same
# number formats / styles the surrounding cells already use.
# ---------------------------------------------------------------------------

# Write a numeric value into $addr, copying the number format from $fmtFrom
# (another cell already carrying the desired style) if the cell's format
# needs to change (e.g. it used to hold a text placeholder).
function Set-Num($addr, $value, $fmtFrom) {
    $dst = $ws.Range($addr)
    $dst.Value = $value
    if ($fmtFrom) {
        $ws.Range($fmtFrom).Copy()
        $dst.PasteSpecial(-4122) | Out-Null
    }
}

# Write one of the placeholder text values ("0" / "***.*") into $addr,
# copying both the *value* and *format* from a donor cell ($donor) that
# already holds that exact placeholder with the correct style. Using
# PasteSpecial values (rather than Range.Value = "0") keeps the cell a true
# shared-string cell instead of letting Excel reinterpret "0" as a number.
function Set-Placeholder($addr, $donor) {
    $dst = $ws.Range($addr)
    $src = $ws.Range($donor)
    $src.Copy()
    $dst.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $src.Copy()
    $dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 15 (Rape): C15, G15 -> "0" placeholder ; H15 -> "***.*" placeholder
# ---------------------------------------------------------------------------
Set-Placeholder "C15" "D15"
Set-Placeholder "G15" "D15"
Set-Placeholder "H15" "E15"

# ---------------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------------
Set-Num "C16" 2
Set-Placeholder "D16" "D15"
Set-Placeholder "E16" "E15"
Set-Num "F16" 10
Set-Num "G16" 4
Set-Num "H16" 150
Set-Num "I16" 7
Set-Num "K16" 75
Set-Num "L16" 16.666666666666
Set-Num "M16" 40
Set-Num "N16" -53.333333333333

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault) -- D17/E17 go from text placeholders back to numbers
# ---------------------------------------------------------------------------
Set-Num "C17" 1
Set-Num "D17" 1 "F16"
Set-Num "E17" 0 "H16"
Set-Num "F17" 5
Set-Num "H17" -16.666666666666
Set-Num "I17" 5
Set-Num "J17" 2
Set-Num "K17" 150
Set-Num "L17" 25
Set-Num "M17" 25
Set-Num "N17" -44.444444444444

# ---------------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------------
Set-Num "C18" 3
Set-Num "D18" 9
Set-Num "E18" -66.666666666666
Set-Num "F18" 13
Set-Num "G18" 26
Set-Num "H18" -50
Set-Num "I18" 11
Set-Num "J18" 21
Set-Num "K18" -47.619047619047
Set-Num "L18" -31.25
Set-Num "M18" -31.25
Set-Num "N18" -81.355932203389

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
Set-Num "C19" 7
Set-Num "E19" -46.153846153846
Set-Num "F19" 24
Set-Num "G19" 56
Set-Num "H19" -57.142857142857
Set-Num "I19" 20
Set-Num "J19" 40
Set-Num "K19" -50
Set-Num "L19" -47.368421052631
Set-Num "M19" -16.666666666666
Set-Num "N19" -25.925925925925

# ---------------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------------
Set-Num "C20" 3
Set-Num "D20" 4
Set-Num "E20" -25
Set-Num "F20" 24
Set-Num "H20" 84.615384615384
Set-Num "I20" 22
Set-Num "J20" 10
Set-Num "K20" 120
Set-Num "L20" 266.666666666667
Set-Num "M20" 340
Set-Num "N20" -89.423076923076

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
Set-Num "C21" 16
Set-Num "D21" 27
Set-Num "E21" -40.740740740740
Set-Num "F21" 77
Set-Num "G21" 105
Set-Num "H21" -26.666666666666
Set-Num "I21" 66
Set-Num "J21" 77
Set-Num "K21" -14.285714285714
Set-Num "L21" -5.714285714285
Set-Num "M21" 22.222222222222
Set-Num "N21" -79.310344827586

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
Set-Num "C24" 9
Set-Num "D24" 8
Set-Num "E24" 12.5
Set-Num "F24" 51
Set-Num "G24" 39
Set-Num "H24" 30.769230769230
Set-Num "I24" 44
Set-Num "J24" 27
Set-Num "K24" 62.962962962963
Set-Num "L24" -8.333333333333
Set-Num "M24" 41.935483870967

# ---------------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------------
Set-Num "C25" 2
Set-Num "D25" 10
Set-Num "E25" -80
Set-Num "F25" 7
Set-Num "G25" 18
Set-Num "H25" -61.111111111111
Set-Num "I25" 4
Set-Num "J25" 16
Set-Num "K25" -75
Set-Num "L25" -66.666666666666

# ---------------------------------------------------------------------------
# Row 26 (UCR Rape*): C26, G26 -> "0" placeholder ; H26 -> "***.*" placeholder
# ---------------------------------------------------------------------------
Set-Placeholder "C26" "D26"
Set-Placeholder "G26" "D26"
Set-Placeholder "H26" "E26"

# ---------------------------------------------------------------------------
# Row 27 (Other Sex Crimes)
# ---------------------------------------------------------------------------
Set-Num "D27" 1 "F16"
Set-Num "E27" -100 "H16"
Set-Placeholder "F27" "D26"
Set-Num "G27" 1 "F16"
Set-Num "H27" -100 "H16"
Set-Num "J27" 1 "F16"
Set-Num "K27" -100 "H16"
